$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.673.14"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.631.57"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("D12").Value = "1.858.00"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "1.603.57"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "26.637.04"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.89%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").Value = "  +4.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").Value = "1.216.01"
$ws.Range("E36").Value = "  +4.79%  "
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.805"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.500"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.793"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "1.766.91"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0511"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.408"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.41%  "
